# finaliza jornada QH 16-12-2014
# Rebuilds "Worksheet" with a title block (title + three subtitle lines),
# a blue header row, and three data rows (2012-2014), adding an
# "Total Impo CIF" column and recomputing "Balanza" as Impo - Expo.

function HexColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate so no stale formatting leaks into the new layout.
$ws.Cells.Clear()

# ---- Title block (rows 2-5), each merged across A:D -----------------------
$title = $ws.Range("A2:D2")
$title.Merge()
$ws.Range("A2").Value = "Balanza comercial"
$title.Font.Size = 12
$title.Font.Name = "Calibri"
$title.Font.Color = HexColor("000000")
$title.Interior.Pattern = 1
$title.Interior.Color = HexColor("F8F5F0")
$title.Interior.PatternColor = HexColor("000000")
$title.Borders.LineStyle = 1
$title.Borders.Color = HexColor("DFD7CA")

$sub1 = $ws.Range("A3:D3")
$sub1.Merge()
$ws.Range("A3").Value = "PAIS ORIGEN: MEXICO"

$sub2 = $ws.Range("A4:D4")
$sub2.Merge()
$ws.Range("A4").Value = "POSICION ARANCELARIA: 21 PREPARACIONES ALIMENTICIAS DIVERSAS"

$sub3 = $ws.Range("A5:D5")
$sub3.Merge()
$ws.Range("A5").Value = "PERIODO: 2012 - 2014"

$subs = $ws.Range("A3:D5")
$subs.Font.Size = 8
$subs.Font.Name = "Calibri"
$subs.Font.Color = HexColor("000000")

# ---- Header row (row 7) -----------------------------------------------
$header = $ws.Range("A7:D7")
$ws.Range("A7").Value = "Periodo"
$ws.Range("B7").Value = "Total Impo CIF (US$)"
$ws.Range("C7").Value = "Total Expo FOB (US$)"
$ws.Range("D7").Value = "Balanza (US$)"

$header.Font.Name = "Calibri"
$header.Font.Size = 11
$header.Font.Bold = $true
$header.Font.Color = HexColor("FFFFFF")
$header.Interior.Pattern = 1
$header.Interior.Color = HexColor("1F497D")
$header.Interior.PatternColor = HexColor("000000")
$header.Borders.Item(8).LineStyle = 1
$header.Borders.Item(8).Color = HexColor("000000")
$header.HorizontalAlignment = -4108

# ---- Data rows (8-10) ---------------------------------------------------
$ws.Range("A8").Value = 2012
$ws.Range("B8").Value = 15518736.32
$ws.Range("C8").Value = 14447812.13
$ws.Range("D8").Value = -1070924.19

$ws.Range("A9").Value = 2013
$ws.Range("B9").Value = 19015241.27
$ws.Range("C9").Value = 13943063.81
$ws.Range("D9").Value = -5072177.46

$ws.Range("A10").Value = 2014
$ws.Range("B10").Value = 16787467.34
$ws.Range("C10").Value = 12292793.74
$ws.Range("D10").Value = -4494673.6

# ---- Page / view tweaks -------------------------------------------------
$ws.PageSetup.Orientation = 1
$ws.Range("D7").Select()
